$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.038.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.92%  '

$ws.Range("D3").Value = "'1.832.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.41%  '

$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = "'239.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.25%  '

$ws.Range("D6").Value = "'0.6710"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.42%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = "'0.07416"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.73%  '

$ws.Range("D9").Value = "'0.2949"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.62%  '

$ws.Range("D10").Value = "'22.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.62%  '

$ws.Range("D11").Value = "'0.07648"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.68%  '

$ws.Range("D12").Value = "'1.834.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.18%  '

$ws.Range("D13").Value = "'5.000"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.94%  '

$ws.Range("D14").Value = "'0.6720"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.97%  '

$ws.Range("D15").Value = "'86.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.69%  '

$ws.Range("D16").Value = "'6.119"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.90%  '

$ws.Range("D17").Value = "'29.048.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.48%  '

$ws.Range("D18").Value = "'0.000008215"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.10%  '

$ws.Range("D19").Value = "'226.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.74%  '

$ws.Range("D20").Value = "'12.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.69%  '

$ws.Range("D21").Value = "'0.9990"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").Value = "'7.299"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.17%  '

$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").Value = "'160.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.23%  '

$ws.Range("D25").Value = "'0.1428"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.85%  '

$ws.Range("D26").Value = "'8.661"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.06%  '

$ws.Range("D27").Value = "'17.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.88%  '

$ws.Range("D28").Value = "'1.499"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.38%  '

$ws.Range("D29").Value = "'4.232"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.54%  '

$ws.Range("D30").Value = "'4.109"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.73%  '

$ws.Range("E31").Value = '  -0.33%  '

$ws.Range("D32").Value = "'0.05379"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.35%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = "'0.7484"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.92%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = "'1.853"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.24%  '

$ws.Range("D35").Value = "'1.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.73%  '

$ws.Range("D36").Value = "'2.683"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").Value = "'1.292.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.92%  '

$ws.Range("D38").Value = "'0.01804"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.80%  '

$ws.Range("D39").Value = "'2.706"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.88%  '

$ws.Range("D40").Value = "'0.9280"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.54%  '

$ws.Range("D41").Value = "'6.075"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.57%  '

$ws.Range("D42").Value = "'0.00000000133"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.37%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = "'104.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.43%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = "'0.9990"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("B45").Value = 'XinFinNetwork'
$ws.Range("C45").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D45").Value = "'0.08208"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +26.00%  '

$ws.Range("D46").Value = "'1.974.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.27%  '

$ws.Range("D47").Value = "'0.5176"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.81%  '

$ws.Range("D48").Value = "'9.373"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.23%  '

$ws.Range("D49").Value = "'63.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.16%  '

$ws.Range("D50").Value = "'1.749"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.78%  '

$ws.Range("D51").Value = "'0.05925"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
